$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 21:22"

# Update countries (re-ordering) and refreshed case numbers

# Row 4
$ws.Range("B4").Value = 757636
$ws.Range("C4").Value = 18844
$ws.Range("D4").Value = 69171
$ws.Range("E4").Value = 648242
$ws.Range("G4").Value = 1209
$ws.Range("H4").Value = 40223

# Row 68
$ws.Range("B68").Value = 1565
$ws.Range("C68").Value = 75
$ws.Range("E68").Value = 1335

# Row 96
$ws.Range("A96").Value = "Guinea"
$ws.Range("B96").Value = 579
$ws.Range("C96").Value = 61
$ws.Range("D96").Value = 87
$ws.Range("E96").Value = 487
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 5

# Row 97
$ws.Range("A97").Value = "Burkina Faso"
$ws.Range("B97").Value = 576
$ws.Range("C97").Value = 11
$ws.Range("D97").Value = 338
$ws.Range("E97").Value = 202
$ws.Range("F97").Value = 0
$ws.Range("H97").Value = 36

# Row 98
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 562
$ws.Range("C98").Value = 14
$ws.Range("D98").Value = 314
$ws.Range("E98").Value = 222
$ws.Range("H98").Value = 26

# Row 99
$ws.Range("A99").Value = "Kirguistan"
$ws.Range("B99").Value = 554
$ws.Range("C99").Value = 48
$ws.Range("D99").Value = 133
$ws.Range("E99").Value = 416
$ws.Range("F99").Value = 5
$ws.Range("H99").Value = 5

# Row 100
$ws.Range("A100").Value = "Nigeria"
$ws.Range("B100").Value = 542
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 166
$ws.Range("E100").Value = 357
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 19

# Row 101
$ws.Range("A101").Value = "Bolivia"
$ws.Range("B101").Value = 520
$ws.Range("C101").Value = 27
$ws.Range("D101").Value = 31
$ws.Range("E101").Value = 457
$ws.Range("F101").Value = 3
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 32

# Row 138
$ws.Range("A138").Value = "Birmania"
$ws.Range("B138").Value = 111
$ws.Range("C138").Value = 13
$ws.Range("E138").Value = 99
$ws.Range("H138").Value = 5

# Row 139
$ws.Range("A139").Value = "Gabon"
$ws.Range("B139").Value = 109
$ws.Range("C139").Value = 1
$ws.Range("E139").Value = 101
$ws.Range("H139").Value = 1

# Row 164
$ws.Range("A164").Value = "Eritrea"
$ws.Range("D164").Value = 3
$ws.Range("H164").Value = 0

# Row 165
$ws.Range("A165").Value = "Puerto Rico"
$ws.Range("D165").Value = 1
$ws.Range("H165").Value = 2

# Row 166
$ws.Range("A166").Value = "Siria"
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 5
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 3

# Row 167
$ws.Range("A167").Value = "Mozambique"
$ws.Range("C167").Value = 4
$ws.Range("D167").Value = 8
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

# Row 184
$ws.Range("A184").Value = "Islas Virgenes de los Estados Unidos"

# Row 185
$ws.Range("A185").Value = "Fiyi"

# Row 201
$ws.Range("A201").Value = "Nicaragua"
$ws.Range("B201").Value = 10
$ws.Range("C201").Value = 1
$ws.Range("D201").Value = 6
$ws.Range("E201").Value = 2
$ws.Range("H201").Value = 2

# Row 202
$ws.Range("A202").Value = "Gambia"
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 6
$ws.Range("H202").Value = 1

# Row 211
$ws.Range("A211").Value = "Sudan del Sur"

# Row 212
$ws.Range("A212").Value = "Santo Tome y Principe"
